# Appends a new test-case row (row 29) to the "Test Cases" sheet,
# mirroring the formatting of the preceding row (row 28), and updates
# the sheet's selection to the new cell, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy the formatting (styles/fill/borders/font) of row 28 down into the
# new row 29 so the new cells pick up the same cell styles (s="10"/"12"/"2").
$ws.Range("A28:E28").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's values.
$ws.Range("A29").Value2 = "ProfileCountryTypeaheadOptionsDisplayTest"
$ws.Range("B29").Value2 = "TBD"

$descriptionText = "Verify that  'country' using type ahead options should display while enter min 2 characters"
$descCell = $ws.Range("C29")
$descCell.Value2 = $descriptionText

# Italicize just the closing quote character (matches the rich-text run
# structure used for similar rows, e.g. the "'primary institution'" /
# "'country'" cells above).
$descCell.Characters(22, 1).Font.Italic = $true
$descCell.Characters(23, $descriptionText.Length - 22).Font.Italic = $false

$ws.Range("D29").Value2 = "Y"
$ws.Range("E29").Value2 = "PASS"

# Move the selection/view to the newly added cell, as in the target sheet.
$ws.Range("D29").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
